# Results_extra.xlsx was regenerated with refreshed texture-analysis
# metrics (Contrast/Correlation/Energy/Homogeneity) for each sample row;
# the ground-truth "Actual Condition" (B/N/M) label per row is unchanged.
# Update the recalculated feature columns B:E for data rows 4-53
# (row 3 is the header, which is unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(4, 0.0021786305014003435, 0.9896198038413121, 0.78794275071223474, 0.99891068474929967),
    @(5, 0.0022311851845481586, 0.98650445735398462, 0.83244622944556768, 0.99888440740772599),
    @(6, 0.0025818682520981264, 0.98854737321217423, 0.77198583955688438, 0.99870906587395092),
    @(7, 0.0051847583774553788, 0.98960538826973954, 0.49604923795102202, 0.99740762081127232),
    @(8, 0.003214435529259103, 0.99357002138898376, 0.49688202118643354, 0.99839278223537042),
    @(9, 0.0031198370995930355, 0.99374838575493196, 0.49784485330895134, 0.99844008145020346),
    @(10, 0.0035288080884523982, 0.99289088289635763, 0.50010581899350126, 0.99823559595577382),
    @(11, 0.0037925370438850717, 0.99239972028556156, 0.4972221796520998, 0.99810373147805753),
    @(12, 0.0027165993488952526, 0.99279807124818298, 0.62008636178712873, 0.99864170032555244),
    @(13, 0.0025990679665828659, 0.9914208452660892, 0.69445616723091042, 0.9987004660167087),
    @(14, 0.0025570242200646134, 0.99270427298617803, 0.64696707213978832, 0.99872148788996762),
    @(15, 0.0032096578307911198, 0.99133831984190346, 0.62624230294165883, 0.99839517108460452),
    @(16, 0.0025866459505661096, 0.99312954556652455, 0.62093170975561784, 0.99870667702471683),
    @(17, 0.0024863142827384623, 0.99337485220897348, 0.6222355174262979, 0.99875684285863064),
    @(18, 0.004525435988873696, 0.99056991355941848, 0.51560167066457474, 0.997737282005563),
    @(19, 0.0055191972702142034, 0.9791855543276291, 0.72934939750090655, 0.99724040136489278),
    @(20, 0.0019005684505637207, 0.99074523767648925, 0.7927419135017485, 0.99904971577471802),
    @(21, 0.0024394928377522267, 0.99458577639148715, 0.54699533043154625, 0.99878025358112388),
    @(22, 0.0048684747388748903, 0.9883282068313306, 0.57804067870974551, 0.99756576263056251),
    @(23, 0.0025780460933237399, 0.99425598900778334, 0.54860527854300234, 0.99871097695333799),
    @(24, 0.0026793333008449838, 0.99462189241698218, 0.49913526393747787, 0.9986603333495776),
    @(25, 0.0026678668245218242, 0.99140042713401932, 0.68710676983208174, 0.99866606658773904),
    @(26, 0.0024767588858024955, 0.99415992696350419, 0.57343213037586549, 0.99876162055709872),
    @(27, 0.0023515831859413357, 0.99369180460118178, 0.62487171075039838, 0.99882420840702935),
    @(28, 0.00308448213092996, 0.99323755616560749, 0.54080558655258204, 0.99845775893453514),
    @(29, 0.005512508492359027, 0.98434328702826845, 0.64243193999861903, 0.99724374575382047),
    @(30, 0.0034475872144966837, 0.99182439593584315, 0.57487225149377041, 0.99827620639275161),
    @(31, 0.0030128166539102118, 0.99130133849615998, 0.65064220897602298, 0.99849359167304486),
    @(32, 0.002726154745831219, 0.9920805027246874, 0.65304796904085205, 0.99863692262708437),
    @(33, 0.0050720046936109748, 0.98509761689238107, 0.65460515635313588, 0.9974639976531946),
    @(34, 0.0046104790216037967, 0.98820018517256836, 0.60468608485972386, 0.99769476048919803),
    @(35, 0.0019550342130987292, 0.99044274964904722, 0.79348846554041019, 0.99902248289345064),
    @(36, 0.0031045484644954893, 0.9862223382158557, 0.77157306222292565, 0.99844772576775231),
    @(37, 0.0024710256476409157, 0.99455097832992179, 0.54405447457945755, 0.99876448717617938),
    @(38, 0.0036769167409598778, 0.9913694604617499, 0.57030104347638078, 0.9981615416295202),
    @(39, 0.0044384818767564012, 0.9896212800548545, 0.56792907036750706, 0.99778075906162167),
    @(40, 0.0016377950348246441, 0.99143931378224959, 0.80704902863605887, 0.99918110248258774),
    @(41, 0.0039425567757797441, 0.98956792570473373, 0.61814653973565692, 0.99802872161211009),
    @(42, 0.0016416171935990306, 0.98575564603301224, 0.88311420885550629, 0.99917919140320055),
    @(43, 0.0023640052019580919, 0.9918225137888288, 0.70855455328138373, 0.998817997399021),
    @(44, 0.0024614702507049493, 0.99362832533223211, 0.61123009843123355, 0.99876926487464757),
    @(45, 0.002020010912263301, 0.99166814592097263, 0.75553972363910959, 0.9989899945438685),
    @(46, 0.0025742239345493533, 0.99467223861739806, 0.51426062571578579, 0.99871288803272529),
    @(47, 0.0022827843280023774, 0.99227886168012958, 0.70206856934861162, 0.99885860783599878),
    @(48, 0.0033883437534936922, 0.98678689817872267, 0.7401849593683727, 0.99830582812325319),
    @(49, 0.0065894017270424418, 0.96628510778972188, 0.79800917809100125, 0.99670529913647865),
    @(50, 0.002030521848892864, 0.99395985466012415, 0.66180257636249462, 0.99898473907555352),
    @(51, 0.0029851060027959093, 0.99028370181840097, 0.68979712788366254, 0.99850744699860217),
    @(52, 0.0023257836142142265, 0.99011240125617828, 0.76245733615164402, 0.99883710819289295),
    @(53, 0.0059425013544775153, 0.97371037485414313, 0.76805302168907841, 0.9970287493227612)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 2).Value2 = $r[1]
    $ws.Cells.Item($rowNum, 3).Value2 = $r[2]
    $ws.Cells.Item($rowNum, 4).Value2 = $r[3]
    $ws.Cells.Item($rowNum, 5).Value2 = $r[4]
}